# Generate Report for Handback
# Update the generated/handoff/handback timestamps recorded on each sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file row.
$wsOverview.Range("G2").Value = "2016-08-30 13:11:21"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the first file row.
$wsZhCn.Range("H2").Value = "2016-08-30 13:11:10"
$wsZhCn.Range("K2").Value = "2016-08-30 13:11:33"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the first file row.
$wsDeDe.Range("H2").Value = "2016-08-30 13:11:21"
$wsDeDe.Range("K2").Value = "2016-08-30 13:11:41"
